$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.096.91"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "3.128.06"
$ws.Range("E3").Value = "  +3.21%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.20%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.123.02"
$ws.Range("E8").Value = "  +3.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.480"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.19%  "

$ws.Range("E15").Value = "  -0.74%  "

$ws.Range("D16").Value = "3.646.21"
$ws.Range("E16").Value = "  +3.22%  "

$ws.Range("D17").Value = "67.075.62"
$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").Value = "3.129.48"
$ws.Range("E19").Value = "  +3.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").Value = "0.0₃0989"
$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("E34").Value = "  -3.02%  "

$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.984"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.80"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.15%  "

$ws.Range("E40").Value = "  +1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.311"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "

$ws.Range("D45").Value = "2.844.46"
$ws.Range("E45").Value = "  +4.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0358"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "381.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.14%  "
